# steel_pathway_v0.1_ClearnSystemWide / steel_pathway_v0.1_emission_multifuel
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) emission_system sheet: rewrite the POSCO1/POSCO2 emission-
#    budget rows (2 & 3) to the new step-down schedule and zero
#    out the last 6 years (V:AA) of the HYUNDAI1/HYUNDAI2 rows.
# ---------------------------------------------------------------
$wsEmissionSystem = $wb.Worksheets.Item("emission_system")

$wsEmissionSystem.Range("B2:K2").Value = 1500000
$wsEmissionSystem.Range("L2:U2").Value = 1000000
$wsEmissionSystem.Range("V2:AA2").Value = 0

$wsEmissionSystem.Range("B3:U3").Value = 1000000
$wsEmissionSystem.Range("V3:AA3").Value = 0

$wsEmissionSystem.Range("V4").Value = 0
$wsEmissionSystem.Range("V5").Value = 0

# ---------------------------------------------------------------
# 2) emission sheet: restore V2:AA2 to 10000000 (undo the taper)
# ---------------------------------------------------------------
$wsEmission = $wb.Worksheets.Item("emission")
$wsEmission.Range("V2:AA2").Value = 10000000
$wsEmission.Range("E21").Select() | Out-Null

# ---------------------------------------------------------------
# 3) Selection / active-sheet bookkeeping to match the saved view
#    state: emission_system becomes the active tab/selection,
#    technology loses its selection focus.
# ---------------------------------------------------------------
$wsTechnology = $wb.Worksheets.Item("technology")
$wsTechnology.Range("E5").Select() | Out-Null

$wsEmissionSystem.Activate() | Out-Null
$wsEmissionSystem.Range("V2:AA2").Select() | Out-Null
